$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.786.45'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.31%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.618.08'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.45%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '210.98'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E6").Value = '  -2.15%  '
$ws.Range("E7").Value = '  +0.33%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '22.92'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.18%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("E11").Value = '  +0.12%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.849.19'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.37%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.618.98'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("E15").Value = '  -1.84%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.47'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.47%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '27.791.17'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.30%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '225.87'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("E19").Value = '  -1.10%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.0₃0711'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").Value = '  -0.57%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.96'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.95%  '
$ws.Range("E24").Value = '  +1.12%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '154.83'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E28").Value = '  -1.65%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.29'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.61%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E32").Value = '  -1.35%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.400.62'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.14%  '
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("E36").Value = '  -2.86%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -1.13%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.551'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.69%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.842'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("E42").Value = '  -2.36%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '65.15'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("E44").Value = '  -2.56%  '
$ws.Range("E45").Value = '  -3.40%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.758.54'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("E47").Value = '  -3.07%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '89.51'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.69%  '
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("E51").Value = '  -0.44%  '
